# Apply updated "want to go" counts (F) and "min ticket price" (G) values
# to the 展览 (Exhibition) and 全部类型 (All types) sheets.
# Both sheets contain identical data tables, so the same updates are applied twice.

$wb = $excel.ActiveWorkbook

# Row -> (F value, G value) updates. G is left untouched where the diff shows no change.
$updates = @(
    @{ Row = 2;  F = 11799; G = 40 },
    @{ Row = 3;  F = 11509; G = 109 },
    @{ Row = 4;  F = 611 },
    @{ Row = 6;  F = 1039 },
    @{ Row = 11; F = 10835 },
    @{ Row = 12; F = 4186 },
    @{ Row = 16; F = 2472 },
    @{ Row = 18; F = 58 },
    @{ Row = 20; F = 133 },
    @{ Row = 21; F = 456 },
    @{ Row = 22; F = 11156 },
    @{ Row = 23; F = 10950 }
)

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates) {
        $ws.Range("F" + $u.Row).Value = $u.F
        if ($u.ContainsKey("G")) {
            $ws.Range("G" + $u.Row).Value = $u.G
        }
    }
}
